$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 49 (date 43881): update Semesterferien dates from 21.02/22.02 to 14.02/15.02
$ws.Range("B49").Value = "Semesterferien(14.02.2020 optimierungen als Headder File)" + [char]10 + "             (15.02.2020 optimierungen als Headder File) kein Erfolg"

# Row 47 (date 43874): update Lauftextes text
$ws.Range("B47").Value = "Testen eines Lauftextes (optimierungen von RTC und DHT22)"

# Update sheet view: scroll so row 31 (A31) is the top-left visible cell,
# then select B47 to match the saved selection state
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B47").Select()
